$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: Task/Sprint Review result was scheduled -> now completed
$ws.Range("D18").Value = "Complete"
$ws.Range("E18").Value = "7/6/2020"
$ws.Range("F18").Value = "Scheduled for 7/6/20, 1PM - Completed"

# Row 19: Result text stays the same content, nothing else to change here

# Row 20: Stakeholder Sprint Review meeting completed with feedback
$ws.Range("D20").Value = "Complete"
$ws.Range("E20").Value = "7/6/2020"
$ws.Range("F20").Value = "Richard and Anu met with us at our Sprint Review on 7/6/20, 1PM.  We received valuable feedback on product features that confirmed our backlog."

# Restore the active cell selection to C1 on the sheet (as seen after the edit)
$ws.Range("C1").Select()
